# Update line-flow results for the 380 kV case (rows 2-25, columns B,D,E,F,G,I,K,L,M,N)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.618628007520158
$ws.Range("D2").Value = 0.09218878428082977
$ws.Range("E2").Value = 0.08563459442923715
$ws.Range("F2").Value = 1.960902129220301
$ws.Range("G2").Value = 0.002525037949218901
$ws.Range("I2").Value = 1.457948466860898
$ws.Range("K2").Value = 1.010141107076464
$ws.Range("L2").Value = 0.2759125721007649
$ws.Range("M2").Value = 0.2058299514185826
$ws.Range("N2").Value = 2.769825489529772

$ws.Range("B3").Value = 0.603318653067646
$ws.Range("D3").Value = 0.0927874162145832
$ws.Range("E3").Value = 0.08236719045959973
$ws.Range("F3").Value = 1.936539090764086
$ws.Range("G3").Value = 0.002529608024369839
$ws.Range("I3").Value = 1.4612687311635
$ws.Range("K3").Value = 0.9230254759955585
$ws.Range("L3").Value = 0.2676792952302378
$ws.Range("M3").Value = 0.200292546409468
$ws.Range("N3").Value = 2.785798471458072

$ws.Range("B4").Value = 0.5942337657181866
$ws.Range("D4").Value = 0.09316928904889821
$ws.Range("E4").Value = 0.08033280062050352
$ws.Range("F4").Value = 1.922656360596633
$ws.Range("G4").Value = 0.002532563185254732
$ws.Range("I4").Value = 1.463823668693415
$ws.Range("K4").Value = 0.8700509239259304
$ws.Range("L4").Value = 0.2627850741520774
$ws.Range("M4").Value = 0.1970042163166497
$ws.Range("N4").Value = 2.796345286285018

$ws.Range("B5").Value = 0.5906111230870437
$ws.Range("D5").Value = 0.09332850378856339
$ws.Range("E5").Value = 0.07949656594656318
$ws.Range("F5").Value = 1.917269065996308
$ws.Range("G5").Value = 0.002533805056369165
$ws.Range("I5").Value = 1.464994648901211
$ws.Range("K5").Value = 0.8485919211350677
$ws.Range("L5").Value = 0.2608310945477683
$ws.Range("M5").Value = 0.1956923172035125
$ws.Range("N5").Value = 2.80082897792407

$ws.Range("B6").Value = 0.5900143974127303
$ws.Range("D6").Value = 0.09335515861503829
$ws.Range("E6").Value = 0.07935727133161308
$ws.Range("F6").Value = 1.916390801452636
$ws.Range("G6").Value = 0.002534013543799155
$ws.Range("I6").Value = 1.465196930248894
$ws.Range("K6").Value = 0.8450364088368758
$ws.Range("L6").Value = 0.2605090799722518
$ws.Range("M6").Value = 0.1954761771947773
$ws.Range("N6").Value = 2.801584708240618

$ws.Range("B7").Value = 0.5941845871312097
$ws.Range("D7").Value = 0.09317142170582349
$ws.Range("E7").Value = 0.08032155218761972
$ws.Range("F7").Value = 1.922582613125883
$ws.Range("G7").Value = 0.002532579781137027
$ws.Range("I7").Value = 1.463838935303002
$ws.Range("K7").Value = 0.8697610009573111
$ws.Range("L7").Value = 0.2627585583513934
$ws.Range("M7").Value = 0.1969864096616334
$ws.Range("N7").Value = 2.796405002899021

$ws.Range("B8").Value = 0.6132841025494145
$ws.Range("D8").Value = 0.09239222247498624
$ws.Range("E8").Value = 0.08451377403241622
$ws.Range("F8").Value = 1.952277993061486
$ws.Range("G8").Value = 0.002526582835881479
$ws.Range("I8").Value = 1.458986136900776
$ws.Range("K8").Value = 0.9799962229496089
$ws.Range("L8").Value = 0.2730402853864291
$ws.Range("M8").Value = 0.2038975035787018
$ws.Range("N8").Value = 2.775179443131712

$ws.Range("B9").Value = 0.6532287681117168
$ws.Range("D9").Value = 0.09097773541076926
$ws.Range("E9").Value = 0.09251654099428919
$ws.Range("F9").Value = 2.019083367273211
$ws.Range("G9").Value = 0.00251600042208649
$ws.Range("I9").Value = 1.453567285125089
$ws.Range("K9").Value = 1.200304228289781
$ws.Range("L9").Value = 0.2944838366644262
$ws.Range("M9").Value = 0.218335140214414
$ws.Range("N9").Value = 2.739428505003815

$ws.Range("B10").Value = 0.6840853477629594
$ws.Range("D10").Value = 0.09000764189842236
$ws.Range("E10").Value = 0.09827123740057786
$ws.Range("F10").Value = 2.073445500873348
$ws.Range("G10").Value = 0.002508935568559459
$ws.Range("I10").Value = 1.452086783023951
$ws.Range("K10").Value = 1.364776578270835
$ws.Range("L10").Value = 0.3110261084698038
$ws.Range("M10").Value = 0.2294821760458419
$ws.Range("N10").Value = 2.716750456849283

$ws.Range("B11").Value = 0.6984488811767449
$ws.Range("D11").Value = 0.08958131868310915
$ws.Range("E11").Value = 0.1008638444543415
$ws.Range("F11").Value = 2.099335074379539
$ws.Range("G11").Value = 0.002505874092857751
$ws.Range("I11").Value = 1.451957013929118
$ws.Range("K11").Value = 1.440187090843892
$ws.Range("L11").Value = 0.3187242491797804
$ws.Range("M11").Value = 0.2346705896227022
$ws.Range("N11").Value = 2.707214634723613

$ws.Range("B12").Value = 0.7039347407521177
$ws.Range("D12").Value = 0.08942203572926566
$ws.Range("E12").Value = 0.1018421047267495
$ws.Range("F12").Value = 2.109306433428713
$ws.Range("G12").Value = 0.002504736574118122
$ws.Range("I12").Value = 1.451986103543483
$ws.Range("K12").Value = 1.468829415154005
$ws.Range("L12").Value = 0.3216642960575484
$ws.Range("M12").Value = 0.2366521954038348
$ws.Range("N12").Value = 2.703716087063654

$ws.Range("B13").Value = 0.7027511901367234
$ws.Range("D13").Value = 0.08945624421038545
$ws.Range("E13").Value = 0.1016315727042389
$ws.Range("F13").Value = 2.107151459625655
$ws.Range("G13").Value = 0.002504980591669939
$ws.Range("I13").Value = 1.451976358528071
$ws.Range("K13").Value = 1.462656937453062
$ws.Range("L13").Value = 0.3210299948000852
$ws.Range("M13").Value = 0.2362246716509944
$ws.Range("N13").Value = 2.704464556578358

$ws.Range("B14").Value = 0.6988992717838585
$ws.Range("D14").Value = 0.08956817114470894
$ws.Range("E14").Value = 0.1009443960168284
$ws.Range("F14").Value = 2.100152062088583
$ws.Range("G14").Value = 0.002505780072541562
$ws.Range("I14").Value = 1.451957839035295
$ws.Range("K14").Value = 1.442541782852061
$ws.Range("L14").Value = 0.3189656286960485
$ws.Range("M14").Value = 0.2348332797284769
$ws.Range("N14").Value = 2.706924551928708

$ws.Range("B15").Value = 0.6965459341498672
$ws.Range("D15").Value = 0.08963701053431983
$ws.Range("E15").Value = 0.1005230278868652
$ws.Range("F15").Value = 2.095886569641024
$ws.Range("G15").Value = 0.002506272612430574
$ws.Range("I15").Value = 1.451956684563378
$ws.Range("K15").Value = 1.430231899004752
$ws.Range("L15").Value = 0.3177043928474319
$ws.Range("M15").Value = 0.2339832072871175
$ws.Range("N15").Value = 2.708446022526033

$ws.Range("B16").Value = 0.6831532074710935
$ws.Range("D16").Value = 0.09003580472645023
$ws.Range("E16").Value = 0.09810130731077038
$ws.Range("F16").Value = 2.071776960673347
$ws.Range("G16").Value = 0.002509138699375534
$ws.Range("I16").Value = 1.452106208482768
$ws.Range("K16").Value = 1.359860304460597
$ws.Range("L16").Value = 0.3105265022166463
$ws.Range("M16").Value = 0.229145463373726
$ws.Range("N16").Value = 2.717389370417351

$ws.Range("B17").Value = 0.675020683899163
$ws.Range("D17").Value = 0.09028428989630655
$ws.Range("E17").Value = 0.09660928570790972
$ws.Range("F17").Value = 2.057284113331662
$ws.Range("G17").Value = 0.002510935892752482
$ws.Range("I17").Value = 1.452337226751311
$ws.Range("K17").Value = 1.316841730473925
$ws.Range("L17").Value = 0.3061674406385038
$ws.Range("M17").Value = 0.2262077479346587
$ws.Range("N17").Value = 2.72307589437753

$ws.Range("B18").Value = 0.6703738405760191
$ws.Range("D18").Value = 0.09042862154041842
$ws.Range("E18").Value = 0.09574873173700382
$ws.Range("F18").Value = 2.049057347678158
$ws.Range("G18").Value = 0.002511983938906653
$ws.Range("I18").Value = 1.452521276960887
$ws.Range("K18").Value = 1.29215417283865
$ws.Range("L18").Value = 0.3036765019299281
$ws.Range("M18").Value = 0.2245291204265172
$ws.Range("N18").Value = 2.726420086774539

$ws.Range("B19").Value = 0.6688057924909003
$ws.Range("D19").Value = 0.09047773176807716
$ws.Range("E19").Value = 0.09545694995456344
$ws.Range("F19").Value = 2.046290628211707
$ws.Range("G19").Value = 0.002512341256809323
$ws.Range("I19").Value = 1.452592381168934
$ws.Range("K19").Value = 1.283804907622653
$ws.Range("L19").Value = 0.3028359070613078
$ws.Range("M19").Value = 0.2239626683011124
$ws.Range("N19").Value = 2.727564982019587

$ws.Range("B20").Value = 0.6758832224723221
$ws.Range("D20").Value = 0.09025769232495495
$ws.Range("E20").Value = 0.09676835964211961
$ws.Range("F20").Value = 2.05881560114652
$ws.Range("G20").Value = 0.002510743094474292
$ws.Range("I20").Value = 1.452307337788412
$ws.Range("K20").Value = 1.321415370846637
$ws.Range("L20").Value = 0.3066297850664057
$ws.Range("M20").Value = 0.2265193274850219
$ws.Range("N20").Value = 2.722462950412137

$ws.Range("B21").Value = 0.7000294083458414
$ws.Range("D21").Value = 0.08953523694399923
$ws.Range("E21").Value = 0.1011463305576399
$ws.Range("F21").Value = 2.102203401910216
$ws.Range("G21").Value = 0.002505544654805431
$ws.Range("I21").Value = 1.451961155145717
$ws.Range("K21").Value = 1.448447749434933
$ws.Range("L21").Value = 0.3195713060835885
$ws.Range("M21").Value = 0.2352415078391985
$ws.Range("N21").Value = 2.706198937710056

$ws.Range("B22").Value = 0.7160824182790293
$ws.Range("D22").Value = 0.08907564114734257
$ws.Range("E22").Value = 0.103987223907918
$ws.Range("F22").Value = 2.131536735629055
$ws.Range("G22").Value = 0.002502274162532944
$ws.Range("I22").Value = 1.452190913293833
$ws.Range("K22").Value = 1.531972492013324
$ws.Range("L22").Value = 0.3281746641600165
$ws.Range("M22").Value = 0.2410402476481721
$ws.Range("N22").Value = 2.696225081399859

$ws.Range("B23").Value = 0.7074898214726204
$ws.Range("D23").Value = 0.08931978490824299
$ws.Range("E23").Value = 0.1024728073869419
$ws.Range("F23").Value = 2.115791350727392
$ws.Range("G23").Value = 0.002504008103821522
$ws.Range("I23").Value = 1.452026547018399
$ws.Range("K23").Value = 1.487347526934627
$ws.Range("L23").Value = 0.3235695744657363
$ws.Range("M23").Value = 0.2379363706095887
$ws.Range("N23").Value = 2.701488252946675

$ws.Range("B24").Value = 0.6754931795130119
$ws.Range("D24").Value = 0.09026971249520965
$ws.Range("E24").Value = 0.09669645095841872
$ws.Range("F24").Value = 2.058122887355779
$ws.Range("G24").Value = 0.002510830212519646
$ws.Range("I24").Value = 1.452320690998299
$ws.Range("K24").Value = 1.319347489834854
$ws.Range("L24").Value = 0.3064207119838755
$ws.Range("M24").Value = 0.2263784302712253
$ws.Range("N24").Value = 2.722739828974483

$ws.Range("B25").Value = 0.6421569925027768
$ws.Range("D25").Value = 0.09134824722577406
$ws.Range("E25").Value = 0.0903741076141813
$ws.Range("F25").Value = 2.000087417693393
$ws.Range("G25").Value = 0.002518737990529189
$ws.Range("I25").Value = 1.454594284425156
$ws.Range("K25").Value = 1.140252921974025
$ws.Range("L25").Value = 0.2885449492884362
$ws.Range("M25").Value = 0.2183449574545851
$ws.Range("N25").Value = 2.748470574761782
